$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 546.2
$ws.Range("I41").Value = 582.75
$ws.Range("K41").Value = 582.75
$ws.Range("M41").Value = -142.75

$ws.Range("H42").Value = 151
$ws.Range("I42").Value = 113.75
$ws.Range("K42").Value = 341.25
$ws.Range("M42").Value = -111.25

$ws.Range("H62").Value = 2400
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776

$ws.Range("H65").Value = 2400
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880

$ws.Range("H138").Value = 3895.2683
$ws.Range("I138").Value = 933.8333
$ws.Range("J138").Value = 5120.6895
$ws.Range("K138").Value = 2801.4999
$ws.Range("L138").Value = 15362.0685
$ws.Range("M138").Value = 2338.5001
$ws.Range("N138").Value = -25642.0685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8385.538
$ws.Range("I32").Value = 5942.853
$ws.Range("K32").Value = 5942.853
$ws.Range("M32").Value = -5655.853

$ws.Range("H45").Value = 1396.1428
$ws.Range("I45").Value = 1055.75
$ws.Range("K45").Value = 1055.75
$ws.Range("M45").Value = -678.75

$ws.Range("H74").Value = 1475.4412
$ws.Range("I74").Value = 1135.6666
$ws.Range("K74").Value = 1135.6666
$ws.Range("M74").Value = -261.6666

$ws.Range("H77").Value = 1475.4412
$ws.Range("I77").Value = 1135.6666
$ws.Range("K77").Value = 5678.333000000001
$ws.Range("M77").Value = -1310.333000000001

$ws.Range("H124").Value = 39880.832
$ws.Range("J124").Value = 39880.832
$ws.Range("L124").Value = 39880.832
$ws.Range("N124").Value = -49700.832

$ws.Range("H132").Value = 3032.3333
$ws.Range("I132").Value = 2665.3333
$ws.Range("J132").Value = 3399.3333
$ws.Range("K132").Value = 7995.999899999999
$ws.Range("L132").Value = 10197.9999
$ws.Range("M132").Value = -5465.999899999999
$ws.Range("N132").Value = -15257.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4905.125
$ws.Range("J86").Value = 5451.5
$ws.Range("L86").Value = 5451.5
$ws.Range("N86").Value = -7697.5

$ws.Range("H89").Value = 4905.125
$ws.Range("J89").Value = 5451.5
$ws.Range("L89").Value = 27257.5
$ws.Range("N89").Value = -38489.5

$ws.Range("H105").Value = 7293.8
$ws.Range("I105").Value = 4848.4287
$ws.Range("K105").Value = 4848.4287
$ws.Range("M105").Value = -3101.4287

$ws.Range("H134").Value = 1986.1364
$ws.Range("I134").Value = 1694.579
$ws.Range("K134").Value = 5083.737
$ws.Range("M134").Value = -2548.737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 423.75
$ws.Range("I22").Value = 297.5
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 297.5
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = 52.5
$ws.Range("N22").Value = -1250

$ws.Range("H31").Value = 3878.9524
$ws.Range("I31").Value = 1390.4445
$ws.Range("K31").Value = 1390.4445
$ws.Range("M31").Value = -1095.4445

$ws.Range("H34").Value = 3878.9524
$ws.Range("I34").Value = 1390.4445
$ws.Range("K34").Value = 1390.4445
$ws.Range("M34").Value = -1188.4445

$ws.Range("H62").Value = 34460.152
$ws.Range("I62").Value = 3998.5833
$ws.Range("K62").Value = 3998.5833
$ws.Range("M62").Value = -3374.5833

$ws.Range("H65").Value = 34460.152
$ws.Range("I65").Value = 3998.5833
$ws.Range("K65").Value = 19992.9165
$ws.Range("M65").Value = -16872.9165

$ws.Range("H122").Value = 3559.8572
$ws.Range("I122").Value = 3449.4614
$ws.Range("K122").Value = 10348.3842
$ws.Range("M122").Value = -7898.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 239.8
$ws.Range("I23").Value = 233
$ws.Range("J23").Value = 250
$ws.Range("K23").Value = 699
$ws.Range("L23").Value = 750
$ws.Range("M23").Value = -464
$ws.Range("N23").Value = -1220

$ws.Range("H34").Value = 1846.5555
$ws.Range("I34").Value = 1069.8334
$ws.Range("J34").Value = 3400
$ws.Range("K34").Value = 3209.5002
$ws.Range("L34").Value = 10200
$ws.Range("M34").Value = -3125.5002
$ws.Range("N34").Value = -10368

$ws.Range("H36").Value = 2062.5
$ws.Range("I36").Value = 2062.5
$ws.Range("K36").Value = 6187.5
$ws.Range("M36").Value = -6018.5

$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2706
$ws.Range("N39").ClearContents()

$ws.Range("H51").Value = 1366.6666
$ws.Range("J51").Value = 1950
$ws.Range("L51").Value = 5850
$ws.Range("N51").Value = -6770

$ws.Range("H55").Value = 500027.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H75").Value = 358.625
$ws.Range("I75").Value = 298.33334
$ws.Range("J75").Value = 394.8
$ws.Range("K75").Value = 895.0000200000001
$ws.Range("L75").Value = 1184.4
$ws.Range("M75").Value = 102.9999799999999
$ws.Range("N75").Value = -3180.4

$ws.Range("H78").Value = 358.625
$ws.Range("I78").Value = 298.33334
$ws.Range("J78").Value = 394.8
$ws.Range("K78").Value = 2685.00006
$ws.Range("L78").Value = 3553.2
$ws.Range("M78").Value = 2306.99994
$ws.Range("N78").Value = -13537.2

$ws.Range("H87").Value = 200
$ws.Range("I87").Value = 200
$ws.Range("K87").Value = 600
$ws.Range("M87").Value = 648

$ws.Range("H90").Value = 200
$ws.Range("I90").Value = 200
$ws.Range("K90").Value = 1800
$ws.Range("M90").Value = 4440

$ws.Range("H128").Value = 3979888
$ws.Range("I128").Value = 3979888
$ws.Range("K128").Value = 11939664
$ws.Range("M128").Value = -11934684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2145.9656
$ws.Range("I102").Value = 1452.8667
$ws.Range("J102").Value = 2888.5715
$ws.Range("K102").Value = 1452.8667
$ws.Range("L102").Value = 2888.5715
$ws.Range("M102").Value = 169.1333
$ws.Range("N102").Value = -6132.5715

$ws.Range("H122").Value = 61934.707
$ws.Range("I122").Value = 2709.1
$ws.Range("J122").Value = 146542.72
$ws.Range("K122").Value = 8127.299999999999
$ws.Range("L122").Value = 439628.16
$ws.Range("M122").Value = -5677.299999999999
$ws.Range("N122").Value = -444528.16

$ws.Range("H123").Value = 125000
$ws.Range("J123").Value = 125000
$ws.Range("L123").Value = 125000
$ws.Range("N123").Value = -129900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2083.6
$ws.Range("I40").Value = 2083.6
$ws.Range("K40").Value = 2083.6
$ws.Range("M40").Value = -1947.6

$ws.Range("H61").Value = 4163.1816
$ws.Range("I61").Value = 4163.1816
$ws.Range("K61").Value = 4163.1816
$ws.Range("M61").Value = -3961.1816

$ws.Range("H100").Value = 1883.2307
$ws.Range("I100").Value = 1720.2222
$ws.Range("K100").Value = 1720.2222
$ws.Range("M100").Value = -1179.2222

$ws.Range("H113").Value = 4163.1816
$ws.Range("I113").Value = 4163.1816
$ws.Range("K113").Value = 4163.1816
$ws.Range("M113").Value = -1993.1816

$ws.Range("H127").Value = 53571.332
$ws.Range("J127").Value = 53571.332
$ws.Range("L127").Value = 53571.332
$ws.Range("N127").Value = -63491.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 795
$ws.Range("I107").Value = 591
$ws.Range("J107").Value = 863
$ws.Range("K107").Value = 1773
$ws.Range("L107").Value = 2589
$ws.Range("M107").Value = 147
$ws.Range("N107").Value = -6429

$ws.Range("H113").Value = 472.8889
$ws.Range("I113").Value = 614.8
$ws.Range("K113").Value = 1844.4
$ws.Range("M113").Value = 325.6000000000001

$ws.Range("H126").Value = 4425.25
$ws.Range("I126").Value = 4344.222
$ws.Range("J126").Value = 4668.3335
$ws.Range("K126").Value = 13032.666
$ws.Range("L126").Value = 14005.0005
$ws.Range("M126").Value = -10562.666
$ws.Range("N126").Value = -18945.0005
